# Scheduled runner update: refresh market/profit values across Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3182.8857
$ws.Range("J17").Value = 3182.8857
$ws.Range("L17").Value = 9548.6571
$ws.Range("N17").Value = -9884.6571

$ws.Range("H33").Value = 161.2
$ws.Range("I33").Value = 161.2
$ws.Range("K33").Value = 161.2
$ws.Range("M33").Value = 67.80000000000001

$ws.Range("H39").Value = 200.66667
$ws.Range("I39").Value = 200.66667
$ws.Range("K39").Value = 602.00001
$ws.Range("M39").Value = -306.00001

$ws.Range("H41").Value = 855.2
$ws.Range("I41").Value = 842
$ws.Range("J41").Value = 875
$ws.Range("K41").Value = 842
$ws.Range("L41").Value = 875
$ws.Range("M41").Value = -402
$ws.Range("N41").Value = -1755

$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("K48").Value = 3000
$ws.Range("M48").Value = -2708

$ws.Range("H53").Value = 265.33334
$ws.Range("I53").Value = 203
$ws.Range("J53").Value = 315.2
$ws.Range("K53").Value = 203
$ws.Range("L53").Value = 315.2
$ws.Range("M53").Value = 434
$ws.Range("N53").Value = -1589.2

$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("K56").Value = 3000
$ws.Range("M56").Value = -2466

$ws.Range("H58").Value = 1743
$ws.Range("I58").Value = 748.8
$ws.Range("K58").Value = 2246.4
$ws.Range("M58").Value = -2096.4

$ws.Range("H76").Value = 6445.364
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 6988.778
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 6988.778
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -7618.778

$ws.Range("H79").Value = 6445.364
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 6988.778
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 6988.778
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -9172.778

$ws.Range("H86").Value = 5733.8335
$ws.Range("J86").Value = 5750
$ws.Range("L86").Value = 5750
$ws.Range("N86").Value = -7996

$ws.Range("H89").Value = 5733.8335
$ws.Range("J89").Value = 5750
$ws.Range("L89").Value = 28750
$ws.Range("N89").Value = -39982

$ws.Range("H96").Value = 3173.375
$ws.Range("I96").Value = 2064.5
$ws.Range("J96").Value = 6500
$ws.Range("K96").Value = 6193.5
$ws.Range("L96").Value = 19500
$ws.Range("M96").Value = -4820.5
$ws.Range("N96").Value = -22246

$ws.Range("H106").Value = 40484.2
$ws.Range("I106").Value = 41649.11
$ws.Range("J106").Value = 30000
$ws.Range("K106").Value = 41649.11
$ws.Range("L106").Value = 30000
$ws.Range("M106").Value = -41018.11
$ws.Range("N106").Value = -31262

$ws.Range("H112").Value = 2193.625
$ws.Range("J112").Value = 2193.625
$ws.Range("L112").Value = 6580.875
$ws.Range("N112").Value = -8796.875

$ws.Range("H135").Value = 827.7941
$ws.Range("I135").Value = 627.24
$ws.Range("J135").Value = 1384.8889
$ws.Range("K135").Value = 5645.16
$ws.Range("L135").Value = 12464.0001
$ws.Range("M135").Value = -3110.16
$ws.Range("N135").Value = -17534.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4807.1133
$ws.Range("J32").Value = 11315.917
$ws.Range("L32").Value = 11315.917
$ws.Range("N32").Value = -11889.917

$ws.Range("H61").Value = 1846.238
$ws.Range("I61").Value = 1651.5294
$ws.Range("K61").Value = 1651.5294
$ws.Range("M61").Value = -1439.5294

$ws.Range("H74").Value = 1339.0278
$ws.Range("I74").Value = 877.09375
$ws.Range("K74").Value = 877.09375
$ws.Range("M74").Value = -3.09375

$ws.Range("H77").Value = 1339.0278
$ws.Range("I77").Value = 877.09375
$ws.Range("K77").Value = 4385.46875
$ws.Range("M77").Value = -17.46875

$ws.Range("H110").Value = 7767.091
$ws.Range("I110").Value = 5604.875
$ws.Range("K110").Value = 5604.875
$ws.Range("M110").Value = -3559.875

$ws.Range("H136").Value = 1846.238
$ws.Range("I136").Value = 1651.5294
$ws.Range("K136").Value = 4954.5882
$ws.Range("M136").Value = -2404.5882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 989.8
$ws.Range("I94").Value = 988.2308
$ws.Range("K94").Value = 988.2308
$ws.Range("M94").Value = -537.2308

$ws.Range("H134").Value = 2948.5454
$ws.Range("I134").Value = 2728.8823
$ws.Range("K134").Value = 8186.646900000001
$ws.Range("M134").Value = -5651.646900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 15999.667
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 22499.5
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 22499.5
$ws.Range("M13").Value = -2861
$ws.Range("N13").Value = -22777.5

$ws.Range("H31").Value = 4273.2744
$ws.Range("I31").Value = 2584.9167
$ws.Range("J31").Value = 5774.037
$ws.Range("K31").Value = 2584.9167
$ws.Range("L31").Value = 5774.037
$ws.Range("M31").Value = -2289.9167
$ws.Range("N31").Value = -6364.037

$ws.Range("H34").Value = 4273.2744
$ws.Range("I34").Value = 2584.9167
$ws.Range("J34").Value = 5774.037
$ws.Range("K34").Value = 2584.9167
$ws.Range("L34").Value = 5774.037
$ws.Range("M34").Value = -2382.9167
$ws.Range("N34").Value = -6178.037

$ws.Range("H58").Value = 3780.9048
$ws.Range("I58").Value = 1593.75
$ws.Range("J58").Value = 5126.846
$ws.Range("K58").Value = 1593.75
$ws.Range("L58").Value = 5126.846
$ws.Range("M58").Value = -1390.75
$ws.Range("N58").Value = -5532.846

$ws.Range("H132").Value = 2161.3774
$ws.Range("I132").Value = 1641.9791
$ws.Range("K132").Value = 4925.9373
$ws.Range("M132").Value = -2395.9373

$ws.Range("H136").Value = 3780.9048
$ws.Range("I136").Value = 1593.75
$ws.Range("J136").Value = 5126.846
$ws.Range("K136").Value = 4781.25
$ws.Range("L136").Value = 15380.538
$ws.Range("M136").Value = -2231.25
$ws.Range("N136").Value = -20480.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62520.812
$ws.Range("I2").Value = 90924.73
$ws.Range("J2").Value = 32.2
$ws.Range("K2").Value = 545548.38
$ws.Range("L2").Value = 193.2
$ws.Range("M2").Value = -545435.38
$ws.Range("N2").Value = -419.2

$ws.Range("H36").Value = 2995.9
$ws.Range("I36").Value = 1662.1111
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 4986.3333
$ws.Range("L36").Value = 45000
$ws.Range("M36").Value = -4817.3333
$ws.Range("N36").Value = -45338

$ws.Range("H81").Value = 1044.3334
$ws.Range("J81").Value = 1466.3334
$ws.Range("L81").Value = 4399.0002
$ws.Range("N81").Value = -6645.0002

$ws.Range("H84").Value = 1044.3334
$ws.Range("J84").Value = 1466.3334
$ws.Range("L84").Value = 13197.0006
$ws.Range("N84").Value = -24429.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4480.5884
$ws.Range("I80").Value = 4017.7
$ws.Range("J80").Value = 5141.857
$ws.Range("K80").Value = 4017.7
$ws.Range("L80").Value = 5141.857
$ws.Range("M80").Value = -3019.7
$ws.Range("N80").Value = -7137.857

$ws.Range("H83").Value = 4480.5884
$ws.Range("I83").Value = 4017.7
$ws.Range("J83").Value = 5141.857
$ws.Range("K83").Value = 20088.5
$ws.Range("L83").Value = 25709.285
$ws.Range("M83").Value = -15096.5
$ws.Range("N83").Value = -35693.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9461.538
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H19").Value = 24358.334
$ws.Range("I19").Value = 16537.5
$ws.Range("K19").Value = 16537.5
$ws.Range("M19").Value = -16367.5

$ws.Range("H30").Value = 2504
$ws.Range("I30").Value = 3172
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 3172
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = -3064
$ws.Range("N30").Value = -716

$ws.Range("H43").Value = 1555437.5
$ws.Range("J43").Value = 2376333.2
$ws.Range("L43").Value = 2376333.2
$ws.Range("N43").Value = -2376719.2

$ws.Range("H46").Value = 3365.9167
$ws.Range("J46").Value = 3898.5
$ws.Range("L46").Value = 3898.5
$ws.Range("N46").Value = -4274.5

$ws.Range("H100").Value = 1467.375
$ws.Range("I100").Value = 947.25
$ws.Range("K100").Value = 947.25
$ws.Range("M100").Value = -406.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10202.5
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 10202.5
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
